# Implemented more methods, trying to fix Karma
# Update the "Functions" reference sheet: the String/Number/Boolean/Array
# support markers for IsNullOrUndefined (row 4) and IsEqualTo/IsNotEqualTo
# (row 5) are fleshed out from a plain "X" into descriptive test labels,
# and the lone "X" markers on rows 6-7 (column D) get the same treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Row 4 (IsNullOrUndefined)
$ws.Range("C4").Value = "Y- Test"
$ws.Range("D4").Value = "X - Test"
$ws.Range("E4").Value = "Y- Test"
$ws.Range("F4").Value = "Y-Test"

# Row 5 (IsEqualTo / IsNotEqualTo)
$ws.Range("C5").Value = "Y- Test"
$ws.Range("D5").Value = "Y- Test"
$ws.Range("E5").Value = "Y- Test"
$ws.Range("F5").Value = "Y-Test"

# Row 6 (IsGreaterThan / IsNotGreaterThan)
$ws.Range("D6").Value = "Y- Test"

# Row 7 (IsGreaterOrEqualTo / IsNotGreaterOrEqualTo)
$ws.Range("D7").Value = "Y- Test"

# Restore the view to where the author left it: scrolled back up to the
# top-left of the sheet with D8 as the active selection.
$ws.Activate()
[void]$ws.Range("D8").Select()
